# Atualização dos gráficos 26082020
$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal": update the last monthly summary row (row 14) ---
$wsMensal = $wb.Worksheets.Item("Mensal")
$wsMensal.Range("A14").Value = 44066
$wsMensal.Range("B14").Value = 120.7
$wsMensal.Range("D14").Value = 1.65

# --- Sheet "Diario": append the new daily rows 384-390 ---
$wsDiario = $wb.Worksheets.Item("Diario")

$newRows = @(
    @(384, 44060, 153.78, 118.73, 29.52),
    @(385, 44061, 158.61, 118.73, 33.59),
    @(386, 44062, 177.33, 118.73, 49.36),
    @(387, 44063, 183.36, 118.73, 54.44),
    @(388, 44064, 175.84, 118.73, 48.1),
    @(389, 44065, 159.39, 118.73, 34.24),
    @(390, 44066, 146.09, 118.73, 23.05)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $wsDiario.Cells.Item($r, 1).Value = $row[1]
    $wsDiario.Cells.Item($r, 2).Value = $row[2]
    $wsDiario.Cells.Item($r, 3).Value = $row[3]
    $wsDiario.Cells.Item($r, 4).Value = $row[4]
    # Match the date-formatted style used by the rest of column A
    $wsDiario.Range("A" + ($r - 1)).Copy()
    $wsDiario.Range("A" + $r).PasteSpecial(-4122)
}
